# Publish terminology IG 2.0.0
# - bump Version / Date on the Metadata sheet
# - rename the old "Concepts" sheet to "Properties" and replace its data
#   with the CodeSystem property definitions (status / effectiveDate)
# - add a brand-new "Concepts" sheet (after "Properties") holding the
#   original concept table that used to live on "Concepts"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet: bump Version + Date
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item(1)
$wsMeta.Range("B3").Value = "1.1.1"

# The Date cell holds a plain text value ("2025-09-22"), but a bare
# string assignment gets auto-recognised as a real date by Excel's
# value parser (and picks up a date number format in the process).
# Route it through a text formula first, then flatten the formula back
# down to a literal value with copy / paste-special so the cell keeps
# its original (unformatted, shared-string) shape.
$wsMeta.Range("B8").Formula = "=""2025-09-22"""
$wsMeta.Range("B8").Copy()
$wsMeta.Range("B8").PasteSpecial(-4163)

# ---------------------------------------------------------------------
# 2. Grab the existing "Concepts" sheet (sheetId 2) - this currently
#    holds the Level/Code/Display/Definition concept table.
# ---------------------------------------------------------------------
$wsOldConcepts = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 3. Free up the "Concepts" name by renaming the original sheet to
#    "Properties" first, then add the new "Concepts" sheet right after
#    it and copy the existing concept table into it verbatim (values +
#    styles).
# ---------------------------------------------------------------------
$wsOldConcepts.Name = "Properties"

$wsNewConcepts = $wb.Worksheets.Add($null, $wsOldConcepts)
$wsNewConcepts.Name = "Concepts"
$wsOldConcepts.Range("A1:D6").Copy($wsNewConcepts.Range("A1:D6"))

# ---------------------------------------------------------------------
# 4. Overwrite the original sheet ("Properties") with the CodeSystem
#    property definitions.
# ---------------------------------------------------------------------

$wsOldConcepts.Range("A1").Value = "Code"
$wsOldConcepts.Range("B1").Value = "Uri"
$wsOldConcepts.Range("C1").Value = "Description"
$wsOldConcepts.Range("D1").Value = "Type"

$wsOldConcepts.Range("A2").Value = "status"
$wsOldConcepts.Range("B2").Value = "http://hl7.org/fhir/concept-properties#status"
$wsOldConcepts.Range("C2").Value = "A property that indicates the status of the concept. One of active, experimental, deprecated, or retired."
$wsOldConcepts.Range("D2").Value = "code"

$wsOldConcepts.Range("A3").Value = "effectiveDate"
$wsOldConcepts.Range("B3").Value = "http://hl7.org/fhir/concept-properties#effectiveDate"
$wsOldConcepts.Range("C3").Value = "The date at which the concept status was last changed."
$wsOldConcepts.Range("D3").Value = "dateTime"

# drop the leftover rows 4-6 (old concept data no longer belongs here)
$wsOldConcepts.Range("A4:D6").Clear()
